# Update cryptocurrency price/volume data on sheet1 (applies diff between
# the previous scrape and the latest GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text content is unambiguous (contains letters, '%', multiple
# '.' groups, spaces, etc.) and will naturally round-trip as text.
$textChanges = @(
    @{Cell="D2"; Value="59.727.54"},
    @{Cell="E2"; Value="  -3.66%  "},
    @{Cell="D3"; Value="3.274.02"},
    @{Cell="E3"; Value="  -4.21%  "},
    @{Cell="E4"; Value="  +0.01%  "},
    @{Cell="E5"; Value="  -4.13%  "},
    @{Cell="E6"; Value="  -7.77%  "},
    @{Cell="E7"; Value="  +0.00%  "},
    @{Cell="D8"; Value="3.270.83"},
    @{Cell="E8"; Value="  -4.28%  "},
    @{Cell="E9"; Value="  -4.00%  "},
    @{Cell="E10"; Value="  -2.98%  "},
    @{Cell="E11"; Value="  -5.27%  "},
    @{Cell="E12"; Value="  -3.59%  "},
    @{Cell="D13"; Value="3.835.33"},
    @{Cell="E13"; Value="  -4.18%  "},
    @{Cell="E14"; Value="  -0.33%  "},
    @{Cell="E15"; Value="  -7.41%  "},
    @{Cell="D16"; Value="3.274.26"},
    @{Cell="E17"; Value="  -5.49%  "},
    @{Cell="D18"; Value="59.843.40"},
    @{Cell="E18"; Value="  -3.52%  "},
    @{Cell="E19"; Value="  -6.81%  "},
    @{Cell="E20"; Value="  -6.14%  "},
    @{Cell="E21"; Value="  -5.43%  "},
    @{Cell="E22"; Value="  -3.00%  "},
    @{Cell="E23"; Value="  -3.16%  "},
    @{Cell="E24"; Value="  -0.09%  "},
    @{Cell="E25"; Value="  -7.38%  "},
    @{Cell="D26"; Value="3.410.39"},
    @{Cell="E26"; Value="  -4.21%  "},
    @{Cell="E27"; Value="  -9.87%  "},
    @{Cell="E28"; Value="  -5.22%  "},
    @{Cell="E29"; Value="  -1.66%  "},
    @{Cell="E30"; Value="  -7.81%  "},
    @{Cell="E31"; Value="  -0.08%  "},
    @{Cell="E32"; Value="  -5.18%  "},
    @{Cell="E33"; Value="  -6.24%  "},
    @{Cell="E34"; Value="  -3.64%  "},
    @{Cell="E35"; Value="  -8.06%  "},
    @{Cell="B36"; Value="NEARProtocol"},
    @{Cell="C36"; Value="https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"},
    @{Cell="E36"; Value="  -8.95%  "},
    @{Cell="B37"; Value="Monero"},
    @{Cell="C37"; Value="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"},
    @{Cell="E37"; Value="  -1.64%  "},
    @{Cell="E38"; Value="  -6.13%  "},
    @{Cell="E39"; Value="  -5.24%  "},
    @{Cell="D40"; Value="3.304.15"},
    @{Cell="E40"; Value="  -4.29%  "},
    @{Cell="E41"; Value="  -7.84%  "},
    @{Cell="E42"; Value="  -16.28%  "},
    @{Cell="E43"; Value="  -2.29%  "},
    @{Cell="E44"; Value="  -4.79%  "},
    @{Cell="E45"; Value="  -7.25%  "},
    @{Cell="E46"; Value="  -5.91%  "},
    @{Cell="E47"; Value="  -7.34%  "},
    @{Cell="E48"; Value="  -0.05%  "},
    @{Cell="D49"; Value="2.314.50"},
    @{Cell="E49"; Value="  -9.35%  "},
    @{Cell="E50"; Value="  -8.15%  "},
    @{Cell="E51"; Value="  -6.80%  "}
)

foreach ($ch in $textChanges) {
    $ws.Range($ch.Cell).Value = $ch.Value
}

# Cells whose new text looks like a plain decimal number (e.g. "0.999").
# Force the cell to text format first so Excel keeps it as a string instead
# of silently converting it to a numeric value, then restore the default
# "Normal" style so no visible formatting change is introduced.
$forceTextChanges = @(
    @{Cell="D4"; Value="0.999"},
    @{Cell="D5"; Value="554.21"},
    @{Cell="D6"; Value="141.06"},
    @{Cell="D12"; Value="0.403"},
    @{Cell="D15"; Value="26.62"},
    @{Cell="D17"; Value="0.0000162"},
    @{Cell="D19"; Value="6.06"},
    @{Cell="D20"; Value="13.66"},
    @{Cell="D21"; Value="8.47"},
    @{Cell="D22"; Value="371.37"},
    @{Cell="D23"; Value="72.89"},
    @{Cell="D24"; Value="1.00"},
    @{Cell="D29"; Value="0.982"},
    @{Cell="D30"; Value="7.05"},
    @{Cell="D33"; Value="7.43"},
    @{Cell="D34"; Value="22.42"},
    @{Cell="D36"; Value="5.03"},
    @{Cell="D37"; Value="165.34"},
    @{Cell="D38"; Value="1.51"},
    @{Cell="D39"; Value="6.59"},
    @{Cell="D41"; Value="0.0723"},
    @{Cell="D42"; Value="25.86"},
    @{Cell="D43"; Value="41.67"},
    @{Cell="D44"; Value="0.743"},
    @{Cell="D50"; Value="6.31"},
    @{Cell="D51"; Value="21.10"}
)

foreach ($ch in $forceTextChanges) {
    $cell = $ws.Range($ch.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $ch.Value
    $cell.Style = "Normal"
}
